# Add two new coding-question rows (39 and 40) to the tracker sheet,
# and update the current view/selection, matching the upstream commit
# "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 39 : "Remove Element" (Leetcode)
# -----------------------------------------------------------------
$ws.Range("A39").Value = 37

# Date column - copy the date-formatted style from the row above,
# then set the actual date value (12/30/2024 -> serial 45656).
$ws.Range("B38").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$row39Date = Get-Date -Year 2024 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("B39").Value = $row39Date.Date

# Problem statement - copy full cell (keeps wrap style + rich text
# formatting template) from C38, then overwrite the text and re-bold
# the trailing "Leetcode" tag.
$ws.Range("C38").Copy()
$ws.Range("C39").PasteSpecial(-4104)
$ws.Range("C39").Value = "Remove given val from array. Leetcode"
$ws.Range("C39").Characters(30, 8).Font.Bold = $true

# Input / Output columns - reuse the monospace style used elsewhere
# for input/output samples.
$ws.Range("E3").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D39").Value = "nums = [3,2,2,3], val = 3"

$ws.Range("E3").Copy()
$ws.Range("E39").PasteSpecial(-4122)
$ws.Range("E39").Value = "2, nums = [2,2,_,_]"

# Solution complexity + who solved it.
$ws.Range("F38").Copy()
$ws.Range("F39").PasteSpecial(-4104)

$ws.Range("G38").Copy()
$ws.Range("G39").PasteSpecial(-4104)

$ws.Rows("39").RowHeight = 60

# -----------------------------------------------------------------
# Row 40 : "Search Insert Position" (Leetcode)
# -----------------------------------------------------------------
$ws.Range("A40").Value = 38

$ws.Range("B38").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$row40Date = Get-Date -Year 2024 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("B40").Value = $row40Date.Date

$ws.Range("C38").Copy()
$ws.Range("C40").PasteSpecial(-4104)
$ws.Range("C40").Value = "Search Insert Position in  array. )(logn) . Leetcode"
$ws.Range("C40").Characters(45, 8).Font.Bold = $true

$ws.Range("E3").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value = "[1,3,5,6], target = 2"

$ws.Range("E3").Copy()
$ws.Range("E40").PasteSpecial(-4122)
$ws.Range("E40").Value = 1

$ws.Range("F38").Copy()
$ws.Range("F40").PasteSpecial(-4104)
$ws.Range("F40").Value = "O(logn)"

$ws.Range("G38").Copy()
$ws.Range("G40").PasteSpecial(-4104)

$ws.Rows("40").RowHeight = 60

# -----------------------------------------------------------------
# View state - scroll/selection moved down to the newly added rows,
# and the workbook window was resized/maximized.
# -----------------------------------------------------------------
$ws.Range("F41").Select() | Out-Null

$wb.Windows.Item(1).Width = 20490
$wb.Windows.Item(1).Height = 7425

Write-Output "Added rows 39-40 and updated view state."
